# Update existing row 136 (was the last data row) with new match data,
# then append a new row 137 with the match that used to be id 8124823
# (previously stored oddly as a shared-string value in column B, now
# stored as a numeric id) plus its full odds data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 137 needs the same formatting as existing data rows (bold/centered/
# bordered "A" id cell, date-formatted "D" cell). Copy the formatting from
# row 136 (which already matches the rest of the table) down to row 137
# before touching any values.
$ws.Range("A136:AB136").Copy() | Out-Null
$ws.Range("A137:AB137").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- Row 136: updated in place -------------------------------------------
$ws.Cells.Item(136, 1).Value  = 134          # A136
$ws.Cells.Item(136, 2).Value  = 8120939      # B136 (now numeric, not text)
$ws.Cells.Item(136, 3).Value  = "India Super League"   # C136
$ws.Cells.Item(136, 4).Value  = 45405.45833333334      # D136
$ws.Cells.Item(136, 5).Value  = "Odisha FC"             # E136
$ws.Cells.Item(136, 6).Value  = "Mohun Bagan SG"        # F136
$ws.Cells.Item(136, 7).Value  = 2             # G136
$ws.Cells.Item(136, 8).Value  = 1             # H136
$ws.Cells.Item(136, 9).Value  = "H"           # I136
$ws.Cells.Item(136, 10).Value = 3.4           # J136
$ws.Cells.Item(136, 11).Value = 3.5           # K136
$ws.Cells.Item(136, 12).Value = 2             # L136
$ws.Cells.Item(136, 13).Value = 3.1           # M136
$ws.Cells.Item(136, 14).Value = 3.4           # N136
$ws.Cells.Item(136, 15).Value = 2.15          # O136
$ws.Cells.Item(136, 16).Value = 0.25          # P136
$ws.Cells.Item(136, 17).Value = 1.9           # Q136
$ws.Cells.Item(136, 18).Value = 1.95          # R136
$ws.Cells.Item(136, 19).Value = 2.75          # S136
$ws.Cells.Item(136, 20).Value = 1.875         # T136
$ws.Cells.Item(136, 21).Value = 1.975         # U136
$ws.Cells.Item(136, 22).Value = 2.1           # V136
$ws.Cells.Item(136, 23).Value = -1            # W136
$ws.Cells.Item(136, 24).Value = -1            # X136
$ws.Cells.Item(136, 25).Value = 0.8999999999999999  # Y136
$ws.Cells.Item(136, 26).Value = -1            # Z136
$ws.Cells.Item(136, 27).Value = 0.4375        # AA136
$ws.Cells.Item(136, 28).Value = -0.5          # AB136

# --- Row 137: new row added at the end ------------------------------------
$ws.Cells.Item(137, 1).Value  = 135           # A137
$ws.Cells.Item(137, 2).Value  = 8124823       # B137 (numeric now)
$ws.Cells.Item(137, 3).Value  = "India Super League"  # C137
$ws.Cells.Item(137, 4).Value  = 45406.45833333334     # D137
$ws.Cells.Item(137, 5).Value  = "East Bengal Club"    # E137
$ws.Cells.Item(137, 6).Value  = "Chennaiyin FC"       # F137
$ws.Cells.Item(137, 7).Value  = 2             # G137
$ws.Cells.Item(137, 8).Value  = 3             # H137
$ws.Cells.Item(137, 9).Value  = "A"           # I137
$ws.Cells.Item(137, 10).Value = 2.05          # J137
$ws.Cells.Item(137, 11).Value = 3.5           # K137
$ws.Cells.Item(137, 12).Value = 3.2           # L137
$ws.Cells.Item(137, 13).Value = 2.3           # M137
$ws.Cells.Item(137, 14).Value = 3.1           # N137
$ws.Cells.Item(137, 15).Value = 3             # O137
$ws.Cells.Item(137, 16).Value = -0.25         # P137
$ws.Cells.Item(137, 17).Value = 2.025         # Q137
$ws.Cells.Item(137, 18).Value = 1.825         # R137
$ws.Cells.Item(137, 19).Value = 2.5           # S137
$ws.Cells.Item(137, 20).Value = 2             # T137
$ws.Cells.Item(137, 21).Value = 1.85          # U137
$ws.Cells.Item(137, 22).Value = -1            # V137
$ws.Cells.Item(137, 23).Value = -1            # W137
$ws.Cells.Item(137, 24).Value = 2             # X137
$ws.Cells.Item(137, 25).Value = -1            # Y137
$ws.Cells.Item(137, 26).Value = 0.825         # Z137
$ws.Cells.Item(137, 27).Value = 1             # AA137
$ws.Cells.Item(137, 28).Value = -1            # AB137
